$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$row = $used.Row + $used.Rows.Count

$ws.Cells.Item($row, 1).Value = "2025-04-28 22:53:58"
$ws.Cells.Item($row, 2).Value = 288
